$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header casing: Nuevos_comentarios -> Nuevos_Comentarios
$ws.Range("I1").Value = "Nuevos_Comentarios"

# Copy the header format (bold font, border, centered/top alignment) from an
# existing header cell onto the two new header cells before setting values.
$ws.Range("A1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)  # xlPasteFormats

# Add two new header columns
$ws.Range("J1").Value = "Proyecto_Origen"
$ws.Range("K1").Value = "Fichero_Origen"

# Populate new columns for the two data rows
$ws.Range("J2").Value = "CERCANÍAS RENFE"
$ws.Range("K2").Value = "EFFAE CBC SALA CERCANIA.xlsx"

$ws.Range("J3").Value = "CERCANÍAS RENFE"
$ws.Range("K3").Value = "EFFAE CBC SALA CERCANIA.xlsx"
